$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.957.69"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "1.891.98"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.019"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.018"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4680"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3917"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.43"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08020"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.78"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").Value = "1.883.30"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.956"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.099"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("E16").Value = "  +1.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06779"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.67%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.35"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001049"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.12"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.017"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("D22").Value = "27.977.83"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.498"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.96"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.349"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.77%  "
$ws.Range("D26").Value = "2.118.90"
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.42"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.03"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.074"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.434"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.63"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9665"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09498"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.382"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.343"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06117"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02245"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.216"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.083"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5974"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1888"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.31"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.272"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5683"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.12"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.933"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06928"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.57"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.069"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.01%  "
